{"js": "// Mapping of old -> new text, in document paragraph order.\n// (Generated from the commit's unified OOXML diff; the document has\n// exactly one paragraph per list entry: a date heading paragraph followed\n// by 100 table-cell paragraphs, each holding a single \"NN+NN=\" / \"NN-NN=\"\n// arithmetic expression run.)\nconst replacements = [\n  [\"2023-02-12 Sunday\", \"2023-02-13 Monday\"],\n  [\"99-16=\", \"56-17=\"],\n  [\"44+10=\", \"41-12=\"],\n  [\"8+45=\", \"17+27=\"],\n  [\"67+0=\", \"96-4=\"],\n  [\"28-15=\", \"20+6=\"],\n  [\"59+26=\", \"49+20=\"],\n  [\"46-17=\", \"94-0=\"],\n  [\"2+50=\", \"87-40=\"],\n  [\"43+39=\", \"72-15=\"],\n  [\"70+4=\", \"47+46=\"],\n  [\"34+64=\", \"78+14=\"],\n  [\"67-20=\", \"63-63=\"],\n  [\"16+39=\", \"53-26=\"],\n  [\"74-74=\", \"51-21=\"],\n  [\"58-10=\", \"56-48=\"],\n  [\"40+53=\", \"55+33=\"],\n  [\"61-61=\", \"29+39=\"],\n  [\"74-41=\", \"18+79=\"],\n  [\"85-45=\", \"57+6=\"],\n  [\"58-57=\", \"72+20=\"],\n  [\"88-64=\", \"19+12=\"],\n  [\"48+35=\", \"61+32=\"],\n  [\"62+0=\", \"56+1=\"],\n  [\"79+16=\", \"47+17=\"],\n  [\"36+57=\", \"32-29=\"],\n  [\"7+50=\", \"88-13=\"],\n  [\"99-49=\", \"61-30=\"],\n  [\"1+93=\", \"76-45=\"],\n  [\"34+50=\", \"54+10=\"],\n  [\"43-14=\", \"99-9=\"],\n  [\"12+60=\", \"84-69=\"],\n  [\"0+41=\", \"67-30=\"],\n  [\"31+66=\", \"79-21=\"],\n  [\"97-74=\", \"25+39=\"],\n  [\"38+57=\", \"91-20=\"],\n  [\"67-45=\", \"83-54=\"],\n  [\"68-63=\", \"63+31=\"],\n  [\"58+32=\", \"62-50=\"],\n  [\"78+2=\", \"79-65=\"],\n  [\"91-40=\", \"56+32=\"],\n  [\"74-0=\", \"42-40=\"],\n  [\"31+18=\", \"55+44=\"],\n  [\"35-35=\", \"88-25=\"],\n  [\"34-2=\", \"44-2=\"],\n  [\"35+1=\", \"10+38=\"],\n  [\"60+0=\", \"17+13=\"],\n  [\"8+46=\", \"49+1=\"],\n  [\"62+28=\", \"78-41=\"],\n  [\"86-26=\", \"92+2=\"],\n  [\"41+13=\", \"88-18=\"],\n  [\"40-3=\", \"45+30=\"],\n  [\"73-66=\", \"17+66=\"],\n  [\"51+36=\", \"39+9=\"],\n  [\"13+2=\", \"48-21=\"],\n  [\"63+14=\", \"16+1=\"],\n  [\"59+32=\", \"70+4=\"],\n  [\"98-24=\", \"0+75=\"],\n  [\"69-4=\", \"11+25=\"],\n  [\"24-7=\", \"92-58=\"],\n  [\"90-53=\", \"92-71=\"],\n  [\"14+80=\", \"88-80=\"],\n  [\"66+19=\", \"8+48=\"],\n  [\"52-15=\", \"7+45=\"],\n  [\"54-24=\", \"34+23=\"],\n  [\"84+6=\", \"84-47=\"],\n  [\"6+72=\", \"95-38=\"],\n  [\"14+33=\", \"73+16=\"],\n  [\"44+26=\", \"10+36=\"],\n  [\"85-83=\", \"82-27=\"],\n  [\"66-18=\", \"46-11=\"],\n  [\"96-25=\", \"94-57=\"],\n  [\"59+4=\", \"43+42=\"],\n  [\"16+27=\", \"94-66=\"],\n  [\"88-75=\", \"17+3=\"],\n  [\"46-21=\", \"40+5=\"],\n  [\"76-73=\", \"60-19=\"],\n  [\"99-71=\", \"39+13=\"],\n  [\"40+55=\", \"80-45=\"],\n  [\"95-41=\", \"13+46=\"],\n  [\"87-58=\", \"98-13=\"],\n  [\"7+61=\", \"54-47=\"],\n  [\"24+30=\", \"60-23=\"],\n  [\"90-61=\", \"5+39=\"],\n  [\"68-40=\", \"38+32=\"],\n  [\"8+75=\", \"66+0=\"],\n  [\"18+31=\", \"62-24=\"],\n  [\"96-56=\", \"73-25=\"],\n  [\"61+29=\", \"50-6=\"],\n  [\"46+20=\", \"39-5=\"],\n  [\"25+48=\", \"83-8=\"],\n  [\"53+27=\", \"66-4=\"],\n  [\"53+27=\", \"56+1=\"],\n  [\"36-7=\", \"60-43=\"],\n  [\"72+22=\", \"71+13=\"],\n  [\"18+63=\", \"23+33=\"],\n  [\"1+0=\", \"38+41=\"],\n  [\"85-72=\", \"85-38=\"],\n  [\"11+21=\", \"99-30=\"],\n  [\"91-1=\", \"3+88=\"],\n  [\"95-72=\", \"63-34=\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nif (paragraphs.items.length !== replacements.length) {\n  throw new Error(\n    `Expected ${replacements.length} paragraphs, found ${paragraphs.items.length}`\n  );\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const para = paragraphs.items[i];\n  if (para.text !== oldText) {\n    throw new Error(\n      `Paragraph ${i}: expected \"${oldText}\" but found \"${para.text}\"`\n    );\n  }\n  if (oldText !== newText) {\n    para.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Mapping of old -> new text, in document paragraph order.\n# (Generated from the commit's unified OOXML diff; the document has a date\n# heading paragraph followed by a 20-row x 5-column table of arithmetic\n# expressions. In the Word COM object model, $d.Paragraphs enumerates one\n# paragraph per table cell PLUS one extra (empty) paragraph for each row's\n# end-of-row mark, so those empty row-end paragraphs are skipped below.)\n$replacements = @(\n    @(\"2023-02-12 Sunday\", \"2023-02-13 Monday\"),\n    @(\"99-16=\", \"56-17=\"),\n    @(\"44+10=\", \"41-12=\"),\n    @(\"8+45=\", \"17+27=\"),\n    @(\"67+0=\", \"96-4=\"),\n    @(\"28-15=\", \"20+6=\"),\n    @(\"59+26=\", \"49+20=\"),\n    @(\"46-17=\", \"94-0=\"),\n    @(\"2+50=\", \"87-40=\"),\n    @(\"43+39=\", \"72-15=\"),\n    @(\"70+4=\", \"47+46=\"),\n    @(\"34+64=\", \"78+14=\"),\n    @(\"67-20=\", \"63-63=\"),\n    @(\"16+39=\", \"53-26=\"),\n    @(\"74-74=\", \"51-21=\"),\n    @(\"58-10=\", \"56-48=\"),\n    @(\"40+53=\", \"55+33=\"),\n    @(\"61-61=\", \"29+39=\"),\n    @(\"74-41=\", \"18+79=\"),\n    @(\"85-45=\", \"57+6=\"),\n    @(\"58-57=\", \"72+20=\"),\n    @(\"88-64=\", \"19+12=\"),\n    @(\"48+35=\", \"61+32=\"),\n    @(\"62+0=\", \"56+1=\"),\n    @(\"79+16=\", \"47+17=\"),\n    @(\"36+57=\", \"32-29=\"),\n    @(\"7+50=\", \"88-13=\"),\n    @(\"99-49=\", \"61-30=\"),\n    @(\"1+93=\", \"76-45=\"),\n    @(\"34+50=\", \"54+10=\"),\n    @(\"43-14=\", \"99-9=\"),\n    @(\"12+60=\", \"84-69=\"),\n    @(\"0+41=\", \"67-30=\"),\n    @(\"31+66=\", \"79-21=\"),\n    @(\"97-74=\", \"25+39=\"),\n    @(\"38+57=\", \"91-20=\"),\n    @(\"67-45=\", \"83-54=\"),\n    @(\"68-63=\", \"63+31=\"),\n    @(\"58+32=\", \"62-50=\"),\n    @(\"78+2=\", \"79-65=\"),\n    @(\"91-40=\", \"56+32=\"),\n    @(\"74-0=\", \"42-40=\"),\n    @(\"31+18=\", \"55+44=\"),\n    @(\"35-35=\", \"88-25=\"),\n    @(\"34-2=\", \"44-2=\"),\n    @(\"35+1=\", \"10+38=\"),\n    @(\"60+0=\", \"17+13=\"),\n    @(\"8+46=\", \"49+1=\"),\n    @(\"62+28=\", \"78-41=\"),\n    @(\"86-26=\", \"92+2=\"),\n    @(\"41+13=\", \"88-18=\"),\n    @(\"40-3=\", \"45+30=\"),\n    @(\"73-66=\", \"17+66=\"),\n    @(\"51+36=\", \"39+9=\"),\n    @(\"13+2=\", \"48-21=\"),\n    @(\"63+14=\", \"16+1=\"),\n    @(\"59+32=\", \"70+4=\"),\n    @(\"98-24=\", \"0+75=\"),\n    @(\"69-4=\", \"11+25=\"),\n    @(\"24-7=\", \"92-58=\"),\n    @(\"90-53=\", \"92-71=\"),\n    @(\"14+80=\", \"88-80=\"),\n    @(\"66+19=\", \"8+48=\"),\n    @(\"52-15=\", \"7+45=\"),\n    @(\"54-24=\", \"34+23=\"),\n    @(\"84+6=\", \"84-47=\"),\n    @(\"6+72=\", \"95-38=\"),\n    @(\"14+33=\", \"73+16=\"),\n    @(\"44+26=\", \"10+36=\"),\n    @(\"85-83=\", \"82-27=\"),\n    @(\"66-18=\", \"46-11=\"),\n    @(\"96-25=\", \"94-57=\"),\n    @(\"59+4=\", \"43+42=\"),\n    @(\"16+27=\", \"94-66=\"),\n    @(\"88-75=\", \"17+3=\"),\n    @(\"46-21=\", \"40+5=\"),\n    @(\"76-73=\", \"60-19=\"),\n    @(\"99-71=\", \"39+13=\"),\n    @(\"40+55=\", \"80-45=\"),\n    @(\"95-41=\", \"13+46=\"),\n    @(\"87-58=\", \"98-13=\"),\n    @(\"7+61=\", \"54-47=\"),\n    @(\"24+30=\", \"60-23=\"),\n    @(\"90-61=\", \"5+39=\"),\n    @(\"68-40=\", \"38+32=\"),\n    @(\"8+75=\", \"66+0=\"),\n    @(\"18+31=\", \"62-24=\"),\n    @(\"96-56=\", \"73-25=\"),\n    @(\"61+29=\", \"50-6=\"),\n    @(\"46+20=\", \"39-5=\"),\n    @(\"25+48=\", \"83-8=\"),\n    @(\"53+27=\", \"66-4=\"),\n    @(\"53+27=\", \"56+1=\"),\n    @(\"36-7=\", \"60-43=\"),\n    @(\"72+22=\", \"71+13=\"),\n    @(\"18+63=\", \"23+33=\"),\n    @(\"1+0=\", \"38+41=\"),\n    @(\"85-72=\", \"85-38=\"),\n    @(\"11+21=\", \"99-30=\"),\n    @(\"91-1=\", \"3+88=\"),\n    @(\"95-72=\", \"63-34=\"),\n)\n\n$d = $word.ActiveDocument\n\n$idx = 0\n$p = 1\n$total = $d.Paragraphs.Count\nwhile ($p -le $total -and $idx -lt $replacements.Length) {\n    $rng = $d.Paragraphs($p).Range\n    $text = $rng.Text\n    # Strip the trailing paragraph/cell-end marker characters for comparison.\n    $trimmed = $text.TrimEnd([char]13, [char]7)\n\n    if ($trimmed -eq \"\" -and $idx -gt 0) {\n        # Row-end paragraph in a table row -- not part of the data, skip it.\n        $p++\n        continue\n    }\n\n    $pair = $replacements[$idx]\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    if ($trimmed -ne $oldText) {\n        throw \"Paragraph $p`: expected '$oldText' but found '$trimmed'\"\n    }\n\n    if ($oldText -ne $newText) {\n        $rng.Text = $newText\n    }\n\n    $idx++\n    $p++\n}\n\nif ($idx -ne $replacements.Length) {\n    throw \"Only applied $idx of $($replacements.Length) replacements\"\n}\n"}
